$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.2
$ws.Range("I2").Value = 1.8
$ws.Range("J2").Value = 4.33
$ws.Range("K2").Value = 2.4
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("Y2").Value = 15
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 15
$ws.Range("AD2").Value = 7.5
$ws.Range("AG2").Value = 126
$ws.Range("AH2").Value = 10
$ws.Range("AI2").Value = 10
$ws.Range("AK2").Value = 15
$ws.Range("AM2").Value = 21
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 67
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7.5
$ws.Range("AW2").Value = 401
$ws.Range("AX2").Value = 4
$ws.Range("AY2").Value = 9
$ws.Range("AZ2").Value = 17
$ws.Range("BC2").Value = 101

# Row 3
$ws.Range("G3").Value = 2.45
$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 3.25
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.38
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 10
$ws.Range("Z3").Value = 23
$ws.Range("AF3").Value = 81
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 34
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 15
$ws.Range("AU3").Value = 9.5
$ws.Range("AX3").Value = 5
$ws.Range("AY3").Value = 21
$ws.Range("BB3").Value = 126

# Row 4
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 4.33
$ws.Range("K4").Value = 2.3
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.25
$ws.Range("AA4").Value = 13
$ws.Range("AI4").Value = 23
$ws.Range("AT4").Value = 3.25
$ws.Range("AU4").Value = 7.5

# Row 5
$ws.Range("G5").Value = 2.8
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 2.25
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 14.1
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 3.45
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 2.2
$ws.Range("S5").Value = 1.29
$ws.Range("T5").Value = 3.4
$ws.Range("W5").Value = 12
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 11
$ws.Range("AB5").Value = 23
$ws.Range("AC5").Value = 15
$ws.Range("AD5").Value = 7.5
$ws.Range("AE5").Value = 12
$ws.Range("AG5").Value = 126
$ws.Range("AH5").Value = 11
$ws.Range("AK5").Value = 21
$ws.Range("AM5").Value = 21

# Row 6
$ws.Range("G6").Value = 1.57
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.1
$ws.Range("K6").Value = 2.5
$ws.Range("L6").Value = 4.75
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 15
$ws.Range("O6").Value = 1.18
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.6
$ws.Range("R6").Value = 2.3
$ws.Range("S6").Value = 1.25
$ws.Range("T6").Value = 3.75
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("AN6").Value = 4
$ws.Range("AO6").Value = 8
$ws.Range("AP6").Value = 15
$ws.Range("AQ6").Value = 21
$ws.Range("AR6").Value = 41
$ws.Range("AS6").Value = 81
$ws.Range("AT6").Value = 3.75
$ws.Range("AU6").Value = 7.5
$ws.Range("AV6").Value = 41
$ws.Range("AW6").Value = 351
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 26
$ws.Range("BA6").Value = 81
$ws.Range("BB6").Value = 81
$ws.Range("BC6").Value = 126

# Row 7
$ws.Range("G7").Value = 4.33
$ws.Range("I7").Value = 1.73
$ws.Range("J7").Value = 4.33
$ws.Range("Q7").Value = 1.5
$ws.Range("R7").Value = 2.5
$ws.Range("U7").Value = 1.5
$ws.Range("V7").Value = 2.5
$ws.Range("X7").Value = 26
$ws.Range("AA7").Value = 29
$ws.Range("AB7").Value = 29
$ws.Range("AH7").Value = 11
$ws.Range("AI7").Value = 11

# Row 8
$ws.Range("K8").Value = 3.4
$ws.Range("O8").Value = 1.07
$ws.Range("P8").Value = 9
$ws.Range("Q8").Value = 1.29
$ws.Range("R8").Value = 3.6
$ws.Range("S8").Value = 1.17
$ws.Range("T8").Value = 5
$ws.Range("X8").Value = 8
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 26
$ws.Range("AI8").Value = 81
$ws.Range("AK8").Value = 251
$ws.Range("AM8").Value = 81
$ws.Range("AN8").Value = 3.5
$ws.Range("AS8").Value = 81
$ws.Range("AT8").Value = 5
$ws.Range("AZ8").Value = 41
$ws.Range("BC8").Value = 301

# Row 15
$ws.Range("G15").Value = 4.1
$ws.Range("H15").Value = 3.4
$ws.Range("I15").Value = 1.9
$ws.Range("J15").Value = 4.5
$ws.Range("L15").Value = 2.6
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.75
$ws.Range("U15").Value = 1.91
$ws.Range("V15").Value = 1.8
$ws.Range("AL15").Value = 17
$ws.Range("AZ15").Value = 23

# Row 16
$ws.Range("G16").Value = 2.15
$ws.Range("I16").Value = 3.75
$ws.Range("J16").Value = 3
$ws.Range("L16").Value = 4.5
$ws.Range("U16").Value = 2.25
$ws.Range("V16").Value = 1.57
$ws.Range("AA16").Value = 21
$ws.Range("AE16").Value = 21
$ws.Range("AH16").Value = 8
$ws.Range("AI16").Value = 17
$ws.Range("AJ16").Value = 15
$ws.Range("AL16").Value = 41
$ws.Range("AQ16").Value = 41
$ws.Range("AX16").Value = 5.5
$ws.Range("AY16").Value = 23

# Row 17
$ws.Range("G17").Value = 2.4
$ws.Range("I17").Value = 3.5
$ws.Range("J17").Value = 3.4
$ws.Range("L17").Value = 4.33
$ws.Range("M17").Value = 1.17
$ws.Range("N17").Value = 5
$ws.Range("O17").Value = 1.62
$ws.Range("P17").Value = 2.2
$ws.Range("Q17").Value = 3.1
$ws.Range("R17").Value = 1.36
$ws.Range("AF17").Value = 81
$ws.Range("AN17").Value = 4.33
$ws.Range("AV17").Value = 81
$ws.Range("AY17").Value = 21
